$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh re-ordered the 25 data rows (rows 2-26, columns A-T)
# into a different sequence without adding/removing any records - i.e. each
# "after" row is an exact copy of some "before" row, just relocated.
# Map: destination row -> source row (based on the content permutation
# derived from the diff).
$writeMap = @{}
$writeMap[2]  = 6
$writeMap[3]  = 7
$writeMap[4]  = 8
$writeMap[5]  = 24
$writeMap[6]  = 10
$writeMap[7]  = 11
$writeMap[8]  = 12
$writeMap[9]  = 13
$writeMap[10] = 17
$writeMap[11] = 14
$writeMap[12] = 2
$writeMap[13] = 3
$writeMap[14] = 18
$writeMap[15] = 23
$writeMap[16] = 15
$writeMap[17] = 16
$writeMap[18] = 4
$writeMap[19] = 5
$writeMap[20] = 19
$writeMap[21] = 20
$writeMap[22] = 21
$writeMap[23] = 22
$writeMap[24] = 25
$writeMap[25] = 26
$writeMap[26] = 9

$firstCol = 1   # A
$lastCol  = 20  # T

# Snapshot every source row's cell values first (Value2 avoids locale/date
# formatting round-trip issues), since the permutation is a single 25-cycle
# and rows would otherwise clobber each other's data if written in place.
$buffer = @{}
for ($r = 2; $r -le 26; $r++) {
    $rowVals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $buffer[$r] = $rowVals
}

# Now write each destination row from the buffered source row's snapshot.
foreach ($destRow in $writeMap.Keys) {
    $srcRow = $writeMap[$destRow]
    $srcVals = $buffer[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c - $firstCol]
    }
}
